$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 13453444
$ws.Range("J17").Value = 13453444
$ws.Range("L17").Value = 40360332
$ws.Range("N17").Value = -40360668
$ws.Range("H80").Value = 4530.1304
$ws.Range("I80").Value = 4127.273
$ws.Range("J80").Value = 4899.4165
$ws.Range("K80").Value = 12381.819
$ws.Range("L80").Value = 14698.2495
$ws.Range("M80").Value = -11383.819
$ws.Range("N80").Value = -16694.2495
$ws.Range("H83").Value = 4530.1304
$ws.Range("I83").Value = 4127.273
$ws.Range("J83").Value = 4899.4165
$ws.Range("K83").Value = 37145.457
$ws.Range("L83").Value = 44094.7485
$ws.Range("M83").Value = -32153.457
$ws.Range("N83").Value = -54078.7485
$ws.Range("H94").Value = 1259.2
$ws.Range("I94").Value = 1259.2
$ws.Range("K94").Value = 1259.2
$ws.Range("M94").Value = -808.2
$ws.Range("H101").Value = 1003
$ws.Range("I101").Value = 1003
$ws.Range("K101").Value = 3009
$ws.Range("M101").Value = -1387
$ws.Range("H131").Value = 1152.4
$ws.Range("I131").Value = 1128.125
$ws.Range("J131").Value = 1249.5
$ws.Range("K131").Value = 3384.375
$ws.Range("L131").Value = 3748.5
$ws.Range("M131").Value = 1655.625
$ws.Range("N131").Value = -13828.5
$ws.Range("H132").Value = 1650.7046
$ws.Range("I132").Value = 1415.6285
$ws.Range("K132").Value = 4246.8855
$ws.Range("M132").Value = -1716.8855
$ws.Range("H138").Value = 4280915
$ws.Range("J138").Value = 4637239
$ws.Range("L138").Value = 13911717
$ws.Range("N138").Value = -13921997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11370637
$ws.Range("I32").Value = 12053639
$ws.Range("K32").Value = 12053639
$ws.Range("M32").Value = -12053352
$ws.Range("H61").Value = 28575972
$ws.Range("I61").Value = 33337302
$ws.Range("J61").Value = 7999.8
$ws.Range("K61").Value = 33337302
$ws.Range("L61").Value = 7999.8
$ws.Range("M61").Value = -33337090
$ws.Range("N61").Value = -8423.799999999999
$ws.Range("H74").Value = 111237500
$ws.Range("I74").Value = 200225500
$ws.Range("K74").Value = 200225500
$ws.Range("M74").Value = -200224626
$ws.Range("H77").Value = 111237500
$ws.Range("I77").Value = 200225500
$ws.Range("K77").Value = 1001127500
$ws.Range("M77").Value = -1001123132
$ws.Range("H132").Value = 35725700
$ws.Range("I132").Value = 11986.846
$ws.Range("K132").Value = 35960.538
$ws.Range("M132").Value = -33430.538
$ws.Range("H136").Value = 28575972
$ws.Range("I136").Value = 33337302
$ws.Range("J136").Value = 7999.8
$ws.Range("K136").Value = 100011906
$ws.Range("L136").Value = 23999.4
$ws.Range("M136").Value = -100009356
$ws.Range("N136").Value = -29099.4
$ws.Range("H139").Value = 69079.22
$ws.Range("J139").Value = 69079.22
$ws.Range("L139").Value = 69079.22
$ws.Range("N139").Value = -79359.22

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2377.5144
$ws.Range("I134").Value = 2328.9666
$ws.Range("K134").Value = 6986.899800000001
$ws.Range("M134").Value = -4451.899800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3500.238
$ws.Range("I58").Value = 3687.1052
$ws.Range("K58").Value = 3687.1052
$ws.Range("M58").Value = -3484.1052
$ws.Range("H86").Value = 5275.25
$ws.Range("I86").Value = 4412.875
$ws.Range("K86").Value = 4412.875
$ws.Range("M86").Value = -3289.875
$ws.Range("H87").Value = 87000
$ws.Range("I87").Value = 84000
$ws.Range("J87").Value = 90000
$ws.Range("K87").Value = 84000
$ws.Range("L87").Value = 90000
$ws.Range("M87").Value = -82814
$ws.Range("N87").Value = -92372
$ws.Range("H89").Value = 5275.25
$ws.Range("I89").Value = 4412.875
$ws.Range("K89").Value = 22064.375
$ws.Range("M89").Value = -16448.375
$ws.Range("H90").Value = 87000
$ws.Range("I90").Value = 84000
$ws.Range("J90").Value = 90000
$ws.Range("K90").Value = 252000
$ws.Range("L90").Value = 270000
$ws.Range("M90").Value = -246072
$ws.Range("N90").Value = -281856
$ws.Range("H99").Value = 7006.125
$ws.Range("J99").Value = 4487
$ws.Range("L99").Value = 4487
$ws.Range("N99").Value = -7483
$ws.Range("H120").Value = 50078.715
$ws.Range("J120").Value = 52110.4
$ws.Range("L120").Value = 52110.4
$ws.Range("N120").Value = -59368.4
$ws.Range("H122").Value = 3089.5
$ws.Range("I122").Value = 3008.375
$ws.Range("K122").Value = 9025.125
$ws.Range("M122").Value = -6575.125
$ws.Range("H126").Value = 7006.125
$ws.Range("J126").Value = 4487
$ws.Range("L126").Value = 13461
$ws.Range("N126").Value = -18401
$ws.Range("H132").Value = 146319
$ws.Range("I132").Value = 157188.61
$ws.Range("K132").Value = 471565.83
$ws.Range("M132").Value = -469035.83
$ws.Range("H135").Value = 84698.60000000001
$ws.Range("J135").Value = 84698.60000000001
$ws.Range("L135").Value = 84698.60000000001
$ws.Range("N135").Value = -94838.60000000001
$ws.Range("H136").Value = 3500.238
$ws.Range("I136").Value = 3687.1052
$ws.Range("K136").Value = 11061.3156
$ws.Range("M136").Value = -8511.3156
$ws.Range("H141").Value = 433781.47
$ws.Range("J141").Value = 451213.78
$ws.Range("L141").Value = 451213.78
$ws.Range("N141").Value = -461573.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1009.8333
$ws.Range("I3").Value = 1009.8333
$ws.Range("K3").Value = 3029.4999
$ws.Range("M3").Value = -2917.4999
$ws.Range("H26").Value = 769.46155
$ws.Range("I26").Value = 142.85715
$ws.Range("J26").Value = 1500.5
$ws.Range("K26").Value = 428.57145
$ws.Range("L26").Value = 4501.5
$ws.Range("M26").Value = -140.57145
$ws.Range("N26").Value = -5077.5
$ws.Range("H128").Value = 116146
$ws.Range("I128").Value = 116146
$ws.Range("K128").Value = 348438
$ws.Range("M128").Value = -343458
$ws.Range("H134").Value = 7343.107
$ws.Range("I134").Value = 1584.6842
$ws.Range("J134").Value = 19499.777
$ws.Range("K134").Value = 4754.0526
$ws.Range("L134").Value = 58499.33099999999
$ws.Range("M134").Value = 315.9474
$ws.Range("N134").Value = -68639.33099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5581.467
$ws.Range("I80").Value = 5349.364
$ws.Range("J80").Value = 6219.75
$ws.Range("K80").Value = 5349.364
$ws.Range("L80").Value = 6219.75
$ws.Range("M80").Value = -4351.364
$ws.Range("N80").Value = -8215.75
$ws.Range("H83").Value = 5581.467
$ws.Range("I83").Value = 5349.364
$ws.Range("J83").Value = 6219.75
$ws.Range("K83").Value = 26746.82
$ws.Range("L83").Value = 31098.75
$ws.Range("M83").Value = -21754.82
$ws.Range("N83").Value = -41082.75
$ws.Range("H113").Value = 3681.5
$ws.Range("J113").Value = 3987.4
$ws.Range("L113").Value = 3987.4
$ws.Range("N113").Value = -8327.4
$ws.Range("H122").Value = 1877.5
$ws.Range("I122").Value = 1747.6842
$ws.Range("K122").Value = 5243.0526
$ws.Range("M122").Value = -2793.0526
$ws.Range("H132").Value = 4231.6743
$ws.Range("I132").Value = 3853.282
$ws.Range("K132").Value = 11559.846
$ws.Range("M132").Value = -9029.846000000001
$ws.Range("H135").Value = 88850.86
$ws.Range("J135").Value = 88850.86
$ws.Range("L135").Value = 88850.86
$ws.Range("N135").Value = -98990.86

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4169.8125
$ws.Range("I22").Value = 2685.75
$ws.Range("J22").Value = 4664.5
$ws.Range("K22").Value = 2685.75
$ws.Range("L22").Value = 4664.5
$ws.Range("M22").Value = -2390.75
$ws.Range("N22").Value = -5254.5
$ws.Range("H27").Value = 4169.8125
$ws.Range("I27").Value = 2685.75
$ws.Range("J27").Value = 4664.5
$ws.Range("K27").Value = 2685.75
$ws.Range("L27").Value = 4664.5
$ws.Range("M27").Value = -2578.75
$ws.Range("N27").Value = -4878.5
$ws.Range("H55").Value = 1120.8948
$ws.Range("I55").Value = 638.5
$ws.Range("J55").Value = 1471.7273
$ws.Range("K55").Value = 638.5
$ws.Range("L55").Value = 1471.7273
$ws.Range("M55").Value = -465.5
$ws.Range("N55").Value = -1817.7273
$ws.Range("H131").Value = 54064.5
$ws.Range("J131").Value = 54064.5
$ws.Range("L131").Value = 54064.5
$ws.Range("N131").Value = -64144.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 72999.336
$ws.Range("J75").Value = 79999
$ws.Range("L75").Value = 79999
$ws.Range("N75").Value = -81871
$ws.Range("H78").Value = 72999.336
$ws.Range("J78").Value = 79999
$ws.Range("L78").Value = 239997
$ws.Range("N78").Value = -249357
$ws.Range("H86").Value = 107258.8
$ws.Range("I86").Value = 55431.668
$ws.Range("J86").Value = 184999.5
$ws.Range("K86").Value = 55431.668
$ws.Range("L86").Value = 184999.5
$ws.Range("M86").Value = -54308.668
$ws.Range("N86").Value = -187245.5
$ws.Range("H89").Value = 107258.8
$ws.Range("I89").Value = 55431.668
$ws.Range("J89").Value = 184999.5
$ws.Range("K89").Value = 277158.34
$ws.Range("L89").Value = 924997.5
$ws.Range("M89").Value = -271542.34
$ws.Range("N89").Value = -936229.5
$ws.Range("H132").Value = 4765.4053
$ws.Range("I132").Value = 4408.1846
$ws.Range("J132").Value = 6423.9287
$ws.Range("K132").Value = 13224.5538
$ws.Range("L132").Value = 19271.7861
$ws.Range("M132").Value = -10694.5538
$ws.Range("N132").Value = -24331.7861
$ws.Range("H136").Value = 1418.921
$ws.Range("I136").Value = 1384.5161
$ws.Range("J136").Value = 1571.2858
$ws.Range("K136").Value = 4153.5483
$ws.Range("L136").Value = 4713.857400000001
$ws.Range("M136").Value = -1603.5483
$ws.Range("N136").Value = -9813.857400000001
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360
